$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.171.45"
$ws.Range("E2").Value = "  +0.31%  "

$ws.Range("D3").Value = "2.558.58"
$ws.Range("E3").Value = "  +0.23%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "584.16"

$ws.Range("D6").Value = "147.58"
$ws.Range("E6").Value = "  -0.33%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.590"
$ws.Range("E8").Value = "  +1.99%  "

$ws.Range("E9").Value = "  +4.20%  "

$ws.Range("D10").Value = "5.65"
$ws.Range("E10").Value = "  +1.16%  "

$ws.Range("E11").Value = "  +0.32%  "

$ws.Range("E12").Value = "  +1.01%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.50"
$ws.Range("E13").Value = "  +1.55%  "

$ws.Range("D14").Value = "3.015.27"
$ws.Range("E14").Value = "  +0.24%  "

$ws.Range("D15").Value = "63.065.97"
$ws.Range("E15").Value = "  +0.25%  "

$ws.Range("E16").Value = "  +5.61%  "

$ws.Range("D17").Value = "2.558.64"
$ws.Range("E17").Value = "  +3.91%  "

$ws.Range("E18").Value = "  -1.32%  "

$ws.Range("D19").Value = "4.43"
$ws.Range("E19").Value = "  +3.81%  "

$ws.Range("D20").Value = "341.79"
$ws.Range("E20").Value = "  +2.33%  "

$ws.Range("E21").Value = "  +0.66%  "

$ws.Range("E22").Value = "  +0.07%  "

$ws.Range("D23").Value = "66.51"
$ws.Range("E23").Value = "  +2.73%  "

$ws.Range("D24").Value = "2.684.16"
$ws.Range("E24").Value = "  +0.08%  "

$ws.Range("E25").Value = "  +2.34%  "

$ws.Range("D26").Value = "0.171"
$ws.Range("E26").Value = "  +0.64%  "

$ws.Range("D27").Value = "8.18"
$ws.Range("E27").Value = "  +13.18%  "

$ws.Range("D28").Value = "8.54"
$ws.Range("E28").Value = "  +2.05%  "

$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.03%  "

$ws.Range("E30").Value = "  +0.02%  "

$ws.Range("E31").Value = "  +8.03%  "

$ws.Range("D32").Value = "0.0₃0828"
$ws.Range("E32").Value = "  +2.12%  "

$ws.Range("D33").Value = "463.58"
$ws.Range("E33").Value = "  +12.85%  "

$ws.Range("D34").Value = "176.85"
$ws.Range("E34").Value = "  -0.12%  "

$ws.Range("E35").Value = "  +2.83%  "

$ws.Range("E36").Value = "  +2.33%  "

$ws.Range("D37").Value = "19.27"
$ws.Range("E37").Value = "  +2.40%  "

$ws.Range("E38").Value = "  +3.45%  "

$ws.Range("D40").Value = "1.75"
$ws.Range("E40").Value = "  +0.42%  "

$ws.Range("D41").Value = "0.999"

$ws.Range("D42").Value = "151.07"
$ws.Range("E42").Value = "  -0.32%  "

$ws.Range("E43").Value = "  +2.27%  "

$ws.Range("D44").Value = "21.07"

$ws.Range("E45").Value = "  +6.68%  "

$ws.Range("D46").Value = "0.614"
$ws.Range("E46").Value = "  +1.87%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0980"
$ws.Range("E47").Value = "  +2.39%  "

$ws.Range("E48").Value = "  +2.24%  "

$ws.Range("E49").Value = "  +0.52%  "

$ws.Range("E50").Value = "  -2.03%  "

$ws.Range("E51").Value = "  -0.24%  "
